$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1/K1 previously held shared-string text ("r"/"s"); replace with the
# correct numeric data-source values used throughout the rest of the column.
$ws.Range("J1").Value = 0.5
$ws.Range("K1").Value = 0.6

# Column K (rows 2-51) was 0.5 for every row; correct it to 0.6 to match
# the refreshed data source.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = 0.6
}

# Update the view: select K1:K51 (previously J2:K51) with K1 active, and
# scroll so the selection is in view (previously frozen at topLeftCell A36).
$ws.Range("K1:K51").Select()
